$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price (column D) updates -------------------------------------
# Prices are stored as text in this sheet (not numbers), so a leading
# apostrophe is used to force a text/quote-prefixed value instead of
# letting Excel coerce the numeric-looking string into a real number.
$priceUpdates = @{
    "D2"  = "244.88"
    "D3"  = "23.05"
    "D4"  = "5.423"
    "D5"  = "0.05972"
    "D6"  = "3.388"
    "D7"  = "0.8086"
    "D8"  = "0.9271"
    "D9"  = "0.1429"
    "D10" = "0.07429"
    "D11" = "0.03381"
    "D12" = "0.03036"
    "D14" = "3.955"
    "D15" = "0.001588"
    "D16" = "0.04827"
    "D17" = "0.0005943"
    "D19" = "0.004159"
    "D20" = "0.0009825"
    "D22" = "3.660"
    "D23" = "6.449"
    "D24" = "2.187"
    "D40" = "0.03928"
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Value = "'" + $priceUpdates[$addr]
}

# --- Rows 41-43: coin listing rotated by one position ---------------------
# Row 41 used to be KickToken, row 42 BKEXToken, row 43 CEJI.
# Now row 41 is BKEXToken, row 42 is CEJI, row 43 is KickToken.
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1075"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002641"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.006213"
$ws.Range("E43").Value = "42KickTokenKICK"

# --- Row 44 (LocalTraders) and Row 45 (CoinLion) value refresh ------------
$ws.Range("D44").Value = "'0.007139"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"

$ws.Range("D45").Value = "'0.00005186"
